# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) sheet gets three new trailing columns appended to its
# existing data frame: date, legislator_name, legislator_id. Every data row
# gets the filing date (2012-02-29), the legislator's name (陳根德) and
# their id (833).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorDate = "2012-02-29"
$legislatorName = "陳根德"
$legislatorId = 833

# Header row
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# Figure out how many data rows currently exist (below the header) and
# stamp every one of them with the same legislator metadata.
$lastRow = $ws.Cells(1, 1).EntireColumn.Cells.SpecialCells(11).Row
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $legislatorDate
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
